$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow value updates, then
# restore protection afterwards.
$ws.Unprotect()

# Update the confidential/disclosure banner text (date 07-13 -> 07-14).
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# Updated Weight (D) / Percent Change (E) values for rows 2-15.
$ws.Range("D2").Value = 0.05898449774798366
$ws.Range("E2").Value = -0.00137624861265262

$ws.Range("D3").Value = 0.02033099402953807
$ws.Range("E3").Value = 0.002511591962905735

$ws.Range("D4").Value = 0.02777835969414475
$ws.Range("E4").Value = 0.004901960784313708

$ws.Range("D5").Value = 0.02890240389651199
$ws.Range("E5").Value = -0.02079322294955721

$ws.Range("D6").Value = 0.02993086833560281
$ws.Range("E6").Value = 0.001531058617672665

$ws.Range("D7").Value = 0.01799125379700429
$ws.Range("E7").Value = -0.006622516556291314

$ws.Range("D8").Value = 0.01063619461610977
$ws.Range("E8").Value = -0.02455838000861688

$ws.Range("D9").Value = 0.01071148004608778
$ws.Range("E9").Value = 0

$ws.Range("D10").Value = 0.06705575049753849
$ws.Range("E10").Value = 0.01057317751808595

$ws.Range("D11").Value = 0.06713038127160363
$ws.Range("E11").Value = 0.01000555864369113

$ws.Range("D12").Value = 0.1521106106630355
$ws.Range("E12").Value = 0.01136207134003575

$ws.Range("D13").Value = 0.3948544045249816
$ws.Range("E13").Value = 0.004310716441072504

$ws.Range("D14").Value = 0.1135828008798575
$ws.Range("E14").Value = 0.006916426512968199

$ws.Range("D15").Value = 0.9999999999999999
$ws.Range("E15").Value = 0.004767204357389643

# Restore sheet protection (password is not recoverable from the legacy
# hash stored in the original file, so re-apply with the known password).
$ws.Protect("D382")
